$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = "275 TOPS (INT8-Sparse)"
    "C3"  = "275 TOPS (INT8-Sparse)"
    "C4"  = "248 TOPS (INT8-Sparse)"
    "C5"  = "200 TOPS (INT8-Sparse)"
    "C6"  = "157 TOPS (INT8-Sparse)"
    "C7"  = "117 TOPS (INT8-Sparse)"
    "C8"  = "67 TOPS (INT8-Sparse)"
    "C9"  = "67 TOPS (INT8-Sparse)"
    "C10" = "34 TOPS (INT8-Sparse)"
    "C14" = "30 TOPS (INT8-Sparse)"
    "C15" = "32 TOPS (INT8-Sparse)"
    "C16" = "32 TOPS (INT8-Sparse)"
    "C17" = "21 TOPS (INT8-Sparse)"
    "C18" = "21 TOPS (INT8-Sparse)"
    "C19" = "1.26 TFLOPS (FP16-Dense)"
    "C20" = "1.33 TFLOPS (FP16-Dense)"
    "C21" = "1.33 TFLOPS (FP16-Dense)"
    "C22" = "1.33 TFLOPS (FP16-Dense)"
    "C23" = "0.472 TFLOPS (FP16-Dense)"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
